$d = $word.ActiveDocument

# The cell "Nombre de proyecto" contains a typo "TweettMonitor" (a
# duplicated "t") split across two runs around the _GoBack bookmark:
#   <w:r>Twee</w:r><bookmarkStart/><bookmarkEnd/><w:r>ttMonitor</w:r>
# It must become "TweetMonitor", with the bookmark ending up between a
# new, separate "t" run and "Monitor":
#   <w:r>Twee</w:r><w:r>t</w:r><bookmarkStart/><bookmarkEnd/><w:r>Monitor</w:r>

# Step 1: shrink the trailing run from "ttMonitor" down to "Monitor"
# (dropping both leading "t" characters) while the bookmark is still
# anchored right in front of it.
$d.Content.Find.Execute("ttMonitor", $true, $false, $false, $false, $false, $true, 1, $false, "Monitor", 2)

# Step 2: re-insert a single "t" right after "Twee" (the text immediately
# preceding the old bookmark position), restoring the visible text to
# "TweetMonitor". Prefer the _GoBack bookmark Word already leaves there;
# fall back to locating "Twee" directly if it is ever missing.
if ($d.Bookmarks.Exists("_GoBack")) {
    $pos = $d.Bookmarks("_GoBack").Range.Start
} else {
    $findRange = $d.Content.Duplicate
    $findRange.Find.Execute("Twee", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $pos = $findRange.End
}

$insPoint = $d.Range($pos, $pos)
$insPoint.InsertAfter("t")

# Force the freshly inserted "t" to live in its own run instead of being
# silently merged back into the preceding "Twee" run, by toggling a
# formatting property on just that character and restoring it.
$newCharRange = $d.Range($pos, $pos + 1)
$newCharRange.Bold = 1
$newCharRange.Bold = 0
